$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2020" column (Q) mirroring the existing "2019" column (P) ---

# Header year cell (Q4), formatted like the other year headers (copy format from P4)
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

# Row 5 (bold first data row) - new value with right/center aligned bold numeric style
$ws.Range("D6").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 25.3
$ws.Range("Q5").Font.Bold = $true

# Row 6
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 17.8

# Row 7 (text placeholder "-")
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = "-"

# Row 8
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 2

# Row 9
$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)
$ws.Range("Q9").Value = 5.5

# Row 10
$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = 0

$excel.CutCopyMode = $false

# --- Row 1 height shrinks slightly (author resaved with updated autofit) ---
$ws.Rows.Item(1).RowHeight = 38.25

# --- Selection moved to N13 on last save ---
$ws.Range("N13").Select()
